$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D59: "Não" -> "Sim"
$ws.Range("D59").Value = "Sim"

# Pre-fill row 62 formats by copying row 61's formatting down one row,
# so the new row's cells land on the same style indexes as the rest of
# the table (no new styles get introduced).
$ws.Range("A61:E61").Copy()
$ws.Range("A62:E62").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 62 data: Agência 10, Tipo Elétrica, Data 22/07/2023,
# Depois da Reforma Sim, Custo Manutenção 3200
$ws.Cells.Item(62, 1).Value = 10
$ws.Cells.Item(62, 2).Value = "Elétrica"
$ws.Cells.Item(62, 3).Value = [DateTime]::FromOADate(45129)
$ws.Cells.Item(62, 4).Value = "Sim"
$ws.Cells.Item(62, 5).Value = 3200

# Scroll/selection state: window scrolled so row 52 is the top row,
# active cell C55 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C55").Select()
